$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump version/status/date/contact for the new release.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value  = "0.4.0-snapshot-1"            # Version
$meta.Cells.Item(6, 2).Value  = "draft"                        # Status
$meta.Cells.Item(8, 2).Value  = "2024-05-23T12:16:26+00:00"    # Date
$meta.Cells.Item(10, 2).Value = "ANS (https://esante.gouv.fr)" # Contact

# ---------------------------------------------------------------------------
# 2. Elements sheet: the "Mapping: RIM Mapping" column (AK) and the
#    "Mapping: Spécification métier vers l'extension
#    RORHealthcareServiceContactTelecom" column (AL) swap places - AL's
#    content moves to AK and vice versa, for the header row and every data
#    row, including their widths.
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = 14
for ($r = 1; $r -le $lastRow; $r++) {
    $akValue = $elements.Cells.Item($r, 37).Value()
    $alValue = $elements.Cells.Item($r, 38).Value()
    $elements.Cells.Item($r, 37).Value = $alValue
    $elements.Cells.Item($r, 38).Value = $akValue
}

$elements.Columns.Item(37).ColumnWidth = 89.9296875
$elements.Columns.Item(38).ColumnWidth = 24.98046875
